# Update "想去人数" (column F) counts across the sheets that mirror the
# same events: "展览" (exhibitions), "演出" (shows) and "全部类型" (all types,
# the combined view). Values bump up slightly, reflecting refreshed stats.
$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetExhibit.Range("F4").Value  = 1277
$sheetExhibit.Range("F6").Value  = 307
$sheetExhibit.Range("F7").Value  = 1111
$sheetExhibit.Range("F9").Value  = 6945
$sheetExhibit.Range("F10").Value = 80
$sheetExhibit.Range("F13").Value = 7837
$sheetExhibit.Range("F16").Value = 5445
$sheetExhibit.Range("F18").Value = 2315
$sheetExhibit.Range("F19").Value = 978
$sheetExhibit.Range("F20").Value = 4543
$sheetExhibit.Range("F21").Value = 271
$sheetExhibit.Range("F25").Value = 314
$sheetExhibit.Range("F26").Value = 234
$sheetExhibit.Range("F28").Value = 2082
$sheetExhibit.Range("F32").Value = 37
$sheetExhibit.Range("F36").Value = 1414
$sheetExhibit.Range("F39").Value = 2147

$sheetShow = $wb.Worksheets.Item("演出")
$sheetShow.Range("F2").Value = 84
$sheetShow.Range("F4").Value = 33

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F6").Value  = 1277
$sheetAll.Range("F7").Value  = 84
$sheetAll.Range("F9").Value  = 307
$sheetAll.Range("F10").Value = 1111
$sheetAll.Range("F12").Value = 6945
$sheetAll.Range("F13").Value = 80
$sheetAll.Range("F16").Value = 7837
$sheetAll.Range("F19").Value = 5445
$sheetAll.Range("F21").Value = 2315
$sheetAll.Range("F22").Value = 978
$sheetAll.Range("F23").Value = 4543
$sheetAll.Range("F24").Value = 271
$sheetAll.Range("F29").Value = 33
$sheetAll.Range("F30").Value = 314
$sheetAll.Range("F31").Value = 234
$sheetAll.Range("F33").Value = 2082
$sheetAll.Range("F37").Value = 37
$sheetAll.Range("F42").Value = 1414
$sheetAll.Range("F45").Value = 2147
